# Apply the dated worksheet update: refresh the date and regenerate all
# two-digit x two-digit multiplication problems/answers in the table.

$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-10 Wednesday", "2024-07-11 Thursday"),
    @("55×14=770",  "42×17=714"),
    @("40×16=640",  "38×90=3420"),
    @("32×65=2080", "22×30=660"),
    @("89×29=2581", "35×73=2555"),
    @("52×99=5148", "95×46=4370"),
    @("31×22=682",  "19×57=1083"),
    @("67×26=1742", "73×16=1168"),
    @("97×71=6887", "24×88=2112"),
    @("46×16=736",  "14×67=938"),
    @("99×38=3762", "84×82=6888"),
    @("79×16=1264", "25×87=2175"),
    @("35×30=1050", "49×64=3136"),
    @("72×32=2304", "53×80=4240"),
    @("59×95=5605", "83×73=6059"),
    @("77×16=1232", "54×98=5292"),
    @("58×49=2842", "25×35=875"),
    @("99×23=2277", "52×32=1664"),
    @("67×40=2680", "91×37=3367"),
    @("92×92=8464", "12×17=204"),
    @("28×71=1988", "87×31=2697"),
    @("84×21=1764", "66×48=3168"),
    @("49×34=1666", "96×19=1824"),
    @("66×30=1980", "54×17=918"),
    @("72×75=5400", "71×73=5183"),
    @("57×54=3078", "11×28=308")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
